$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'74.156.00"
$ws.Range("E2").Value = "  +7.72%  "
$ws.Range("D3").Value = "'2.633.40"
$ws.Range("E3").Value = "  +7.98%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'186.10"
$ws.Range("E5").Value = "  +14.58%  "
$ws.Range("D6").Value = "'582.42"
$ws.Range("E6").Value = "  +4.11%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +4.68%  "
$ws.Range("E9").Value = "  +17.07%  "
$ws.Range("D10").Value = "'2.633.12"
$ws.Range("E10").Value = "  +8.02%  "
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "'0.357"
$ws.Range("E12").Value = "  +7.88%  "
$ws.Range("D13").Value = "'4.77"
$ws.Range("E13").Value = "  +3.39%  "
$ws.Range("D14").Value = "'0.0000189"
$ws.Range("E14").Value = "  +5.85%  "
$ws.Range("D15").Value = "'74.000.55"
$ws.Range("E15").Value = "  +7.65%  "
$ws.Range("D16").Value = "'3.116.25"
$ws.Range("E16").Value = "  +7.91%  "
$ws.Range("D17").Value = "'26.28"
$ws.Range("E17").Value = "  +12.77%  "
$ws.Range("D18").Value = "'2.630.78"
$ws.Range("E18").Value = "  +7.81%  "
$ws.Range("D19").Value = "'8.99"
$ws.Range("E19").Value = "  +28.93%  "
$ws.Range("E20").Value = "  +11.45%  "
$ws.Range("D21").Value = "'372.27"
$ws.Range("E21").Value = "  +9.68%  "
$ws.Range("E22").Value = "  +18.19%  "
$ws.Range("D23").Value = "'4.08"
$ws.Range("E23").Value = "  +6.46%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "'70.09"
$ws.Range("E25").Value = "  +5.68%  "
$ws.Range("D26").Value = "'4.15"
$ws.Range("E26").Value = "  +11.12%  "
$ws.Range("D27").Value = "'9.35"
$ws.Range("E27").Value = "  +13.40%  "
$ws.Range("D28").Value = "'2.769.20"
$ws.Range("E28").Value = "  +7.87%  "
$ws.Range("D29").Value = "'1.01"
$ws.Range("E29").Value = "  +2.03%  "
$ws.Range("D30").Value = "'0.0₃0949"
$ws.Range("E30").Value = "  +15.25%  "
$ws.Range("D31").Value = "'525.92"
$ws.Range("D32").Value = "'1.39"
$ws.Range("E32").Value = "  +18.47%  "
$ws.Range("D33").Value = "'7.71"
$ws.Range("E33").Value = "  +7.82%  "
$ws.Range("E34").Value = "  +8.67%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "'163.13"
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("E37").Value = "  +12.74%  "
$ws.Range("D38").Value = "'19.15"
$ws.Range("E38").Value = "  +6.24%  "
$ws.Range("D39").Value = "'19.29"
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("D41").Value = "'4.92"
$ws.Range("E41").Value = "  +12.68%  "
$ws.Range("D42").Value = "'0.328"
$ws.Range("E42").Value = "  +9.39%  "
$ws.Range("E43").Value = "  +10.08%  "
$ws.Range("D44").Value = "'160.42"
$ws.Range("E44").Value = "  +23.09%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "'2.39"
$ws.Range("E45").Value = "  +15.08%  "
$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").Value = "'1.19"
$ws.Range("E46").Value = "  +11.45%  "
$ws.Range("D47").Value = "'38.95"
$ws.Range("E47").Value = "  +3.67%  "
$ws.Range("E48").Value = "  +18.80%  "
$ws.Range("E49").Value = "  +8.65%  "
$ws.Range("D50").Value = "'0.529"
$ws.Range("E50").Value = "  +9.62%  "
$ws.Range("D51").Value = "'20.76"
$ws.Range("E51").Value = "  +22.58%  "

# Clear the "quote prefix" text-format marker the apostrophe-entry above
# leaves behind on the D-column cells, so their style stays the plain
# default (matching the original workbook) while the value remains text.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
